$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = 0.0085
$ws.Range("P2").Value = 0.0071
$ws.Range("Q2").Value = 0.0062
$ws.Range("O5").Value = -0.4557
$ws.Range("P5").Value = -0.1629
$ws.Range("Q5").Value = -0.1465
$ws.Range("O6").Value = -0.0256
$ws.Range("P6").Value = -0.0302
$ws.Range("Q6").Value = -0.0229
$ws.Range("R6").Value = -0.0161
$ws.Range("S6").Value = 0.0049
$ws.Range("T6").Value = -0.0031
$ws.Range("U6").Value = -0.0127
$ws.Range("V6").Value = -0.0274
$ws.Range("W6").Value = -1.7933
$ws.Range("O7").Value = -0.1279
$ws.Range("P7").Value = 0.0584
$ws.Range("Q7").Value = 0.0104
$ws.Range("R7").Value = 0.0486
$ws.Range("S7").Value = -0.0036
$ws.Range("T7").Value = 0.0157
$ws.Range("U7").Value = 0.033
$ws.Range("V7").Value = 0.0476
$ws.Range("W7").Value = -5.2081
$ws.Range("O8").Value = 0.1788
$ws.Range("P8").Value = 0.3287
$ws.Range("Q8").Value = 0.2853
$ws.Range("R8").Value = -0.0597
$ws.Range("S8").Value = -0.0771
$ws.Range("T8").Value = -0.046
$ws.Range("U8").Value = -0.0541
$ws.Range("V8").Value = 0.0282
$ws.Range("W8").Value = 7.6011
$ws.Range("S9").Value = -0.014
$ws.Range("S10").Value = -0.01
$ws.Range("U10").Value = -0.0077
$ws.Range("V10").Value = -0.0071
$ws.Range("O11").Value = 0.0095
$ws.Range("P11").Value = 0.0082
$ws.Range("Q11").Value = 0.0077
$ws.Range("O12").Value = 0.0152
$ws.Range("P12").Value = 0.0421
$ws.Range("Q12").Value = 0.0334
$ws.Range("R12").Value = -0.1148
$ws.Range("S12").Value = -0.0988
$ws.Range("T12").Value = -0.0906
$ws.Range("S13").Value = -0.0236
$ws.Range("O14").Value = -0.2931
$ws.Range("P14").Value = -0.2649
$ws.Range("Q14").Value = -0.1078
$ws.Range("R14").Value = -0.0597
$ws.Range("S14").Value = -0.0948
$ws.Range("T14").Value = -0.0707
$ws.Range("U14").Value = -0.0322
$ws.Range("V14").Value = -0.0131
$ws.Range("W14").Value = -0.0742
$ws.Range("O16").Value = -0.5629
$ws.Range("P16").Value = -0.2064
$ws.Range("Q16").Value = 0.094
$ws.Range("R16").Value = -0.4702
$ws.Range("S16").Value = -0.6931
$ws.Range("T16").Value = -0.4406
$ws.Range("U16").Value = -0.3915
$ws.Range("V16").Value = -0.282
$ws.Range("W16").Value = -73.7292
$ws.Range("O18").Value = -0.2316
$ws.Range("P18").Value = -0.0305
$ws.Range("Q18").Value = 0.0015
$ws.Range("O22").Value = 0.0493
$ws.Range("P22").Value = -0.1489
$ws.Range("Q22").Value = -0.1368
$ws.Range("O23").Value = -0.0165
$ws.Range("P23").Value = -0.0088
$ws.Range("Q23").Value = 0.0082
$ws.Range("R23").Value = 0.0227
$ws.Range("S23").Value = 0.0023
$ws.Range("T23").Value = 0.0025
$ws.Range("U23").Value = 0.0078
$ws.Range("V23").Value = 0.0127
$ws.Range("W23").Value = 0.0755
$ws.Range("O24").Value = 0.1152
$ws.Range("P24").Value = 0.0565
$ws.Range("Q24").Value = 0.1805
$ws.Range("R24").Value = 0.0073
$ws.Range("S24").Value = -0.0236
$ws.Range("T24").Value = -0.0165
$ws.Range("U24").Value = -0.0087
$ws.Range("W24").Value = -1.0114
$ws.Range("O25").Value = 0.2563
$ws.Range("P25").Value = 0.175
$ws.Range("Q25").Value = 0.1406
$ws.Range("R25").Value = 0.0987
$ws.Range("S25").Value = 0.0243
$ws.Range("T25").Value = 0.0542
$ws.Range("U25").Value = 0.0646
$ws.Range("V25").Value = 0.0267
$ws.Range("W25").Value = 3.9236
$ws.Range("O26").Value = 0.0088
$ws.Range("P26").Value = 0.008
$ws.Range("Q26").Value = 0.007
$ws.Range("O29").Value = -0.0014
$ws.Range("P29").Value = -0.0003
$ws.Range("Q29").Value = -0.0025
$ws.Range("R29").Value = -0.0034
$ws.Range("S29").Value = -0.0031
$ws.Range("T29").Value = -0.0024
$ws.Range("U29").Value = -0.0033
$ws.Range("V29").Value = -0.0031
$ws.Range("W29").Value = -0.0618
$ws.Range("O30").Value = 0.0017
$ws.Range("P30").Value = 0.0011
$ws.Range("Q30").Value = 0.0007
$ws.Range("O33").Value = -0.0959
$ws.Range("P33").Value = -0.0956
$ws.Range("Q33").Value = -0.0627
$ws.Range("O34").Value = 0.0947
$ws.Range("P34").Value = 0.2079
$ws.Range("Q34").Value = 0.2006
$ws.Range("R34").Value = 0.2678
$ws.Range("S34").Value = 0.2243
$ws.Range("T34").Value = 0.1675
$ws.Range("U34").Value = 0.0227
$ws.Range("V34").Value = 0.0867
$ws.Range("W34").Value = -1.366
$ws.Range("O35").Value = 0.0553
$ws.Range("P35").Value = 0.0191
$ws.Range("Q35").Value = 0.0115
$ws.Range("R35").Value = -0.0065
$ws.Range("S35").Value = -0.0484
$ws.Range("T35").Value = -0.0164
$ws.Range("U35").Value = -0.0155
$ws.Range("V35").Value = -0.0444
$ws.Range("W35").Value = -0.1839
$ws.Range("O36").Value = 0.0688
$ws.Range("P36").Value = 0.0558
$ws.Range("Q36").Value = 0.0767
$ws.Range("R36").Value = 0.0741
$ws.Range("S36").Value = 0.0746
$ws.Range("T36").Value = 0.0781
$ws.Range("U36").Value = 0.0633
$ws.Range("V36").Value = 0.0956
$ws.Range("W36").Value = 0.1382
$ws.Range("V37").Value = -0.0003
$ws.Range("P38").Value = -0.0003
$ws.Range("Q38").Value = -0.0003
$ws.Range("R38").Value = -0.0002
$ws.Range("V38").Value = -0.0002
$ws.Range("O39").Value = 0.0017
$ws.Range("P39").Value = 0.0016
$ws.Range("Q39").Value = 0.0007
$ws.Range("O40").Value = -0.0091
$ws.Range("P40").Value = -0.0241
$ws.Range("Q40").Value = -0.0365
$ws.Range("R40").Value = -0.0549
$ws.Range("S40").Value = -0.0441
$ws.Range("T40").Value = -0.0366
$ws.Range("Q41").Value = -0.0003
$ws.Range("R41").Value = -0.0004
$ws.Range("O42").Value = -0.0011
$ws.Range("P42").Value = -0.0021
$ws.Range("Q42").Value = -0.0017
$ws.Range("R42").Value = -0.0015
$ws.Range("S42").Value = -0.0021
$ws.Range("T42").Value = -0.0019
$ws.Range("U42").Value = -0.0015
$ws.Range("V42").Value = -0.0015
$ws.Range("W42").Value = -0.0018
$ws.Range("O44").Value = 0.1109
$ws.Range("P44").Value = 0.1567
$ws.Range("Q44").Value = 0.213
$ws.Range("R44").Value = 0.2161
$ws.Range("S44").Value = 0.1264
$ws.Range("T44").Value = 0.1236
$ws.Range("U44").Value = -0.021
$ws.Range("V44").Value = -0.1911
$ws.Range("W44").Value = -3.082
$ws.Range("O46").Value = -0.021
$ws.Range("P46").Value = -0.0198
$ws.Range("Q46").Value = -0.0137
$ws.Range("O50").Value = -0.1065
$ws.Range("P50").Value = -0.1052
$ws.Range("Q50").Value = -0.0743
$ws.Range("O51").Value = -0.0004
$ws.Range("P51").Value = -0.0003
$ws.Range("Q51").Value = 0
$ws.Range("R51").Value = 0.0002
$ws.Range("S51").Value = 0.0001
$ws.Range("T51").Value = 0.0002
$ws.Range("U51").Value = 0.0004
$ws.Range("V51").Value = 0.0009
$ws.Range("W51").Value = 0.0026
$ws.Range("O52").Value = 0.0876
$ws.Range("P52").Value = 0.084
$ws.Range("Q52").Value = 0.0833
$ws.Range("R52").Value = 0.0133
$ws.Range("S52").Value = -0.0038
$ws.Range("T52").Value = -0.0012
$ws.Range("U52").Value = -0.0011
$ws.Range("W52").Value = -0.1176
$ws.Range("O53").Value = 0.0251
$ws.Range("P53").Value = 0.0277
$ws.Range("Q53").Value = 0.0298
$ws.Range("R53").Value = 0.0297
$ws.Range("S53").Value = 0.0296
$ws.Range("T53").Value = 0.032
$ws.Range("U53").Value = 0.0339
$ws.Range("V53").Value = 0.0487
$ws.Range("W53").Value = 0.092
$ws.Range("O54").Value = 0.0017
$ws.Range("P54").Value = 0.0016
$ws.Range("Q54").Value = 0.0007
$ws.Range("O57").Value = -0.0058
$ws.Range("P57").Value = -0.0063
$ws.Range("Q57").Value = -0.0082
$ws.Range("R57").Value = -0.0096
$ws.Range("S57").Value = -0.0094
$ws.Range("T57").Value = -0.0073
$ws.Range("U57").Value = -0.0057
$ws.Range("V57").Value = -0.0044
$ws.Range("W57").Value = 0.0094
